# Update the "Förändrad" (Changed) date column C for rows 2-20 from
# 2023-09-16 (serial 45185) to 2023-10-05 (serial 45204).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 20; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45185) {
        $cell.Value2 = 45204
    }
}
